$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) "Activity" column header becomes "Target" (column C)
$ws.Range("C1").Value = "Target"

# 2) Insert a brand-new column before "Suggestion" (currently column E) and
#    give it the header "Severity". This pushes "Suggestion" to column F.
$ws.Range("E1:E2").EntireColumn.Insert()
$ws.Range("E1").Value = "Severity"

# 3) Grow the table so it covers the newly inserted column too.
$lo.Resize($ws.Range("A1:F2"))

# Re-affirm the "Suggestion" header text on its new location (F1) so the
# table definition picks up the right column name instead of a default.
$ws.Range("F1").Value = "Suggestion"

# Cosmetic: match the column widths from the authored template as closely
# as this engine's character-width rounding allows.
$ws.Columns.Item(1).ColumnWidth = 18.857142857142858
$ws.Columns.Item(3).ColumnWidth = 19.714285714285715
$ws.Columns.Item(4).ColumnWidth = 35.42857142857143
$ws.Columns.Item(5).ColumnWidth = 19.857142857142858
$ws.Columns.Item(6).ColumnWidth = 50.42857142857143

# Cosmetic: the saved workbook had the second row selected.
$ws.Range("A2").Select() | Out-Null
